$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row changes -----------------------------------------------
# Column J ("Mission") now holds "Revenue"
$ws.Range("J1").Value = "Revenue"

# New column K holds the old "Mission" header, renamed "Mission Statement"
$ws.Range("K1").Value = "Mission Statement"

# --- Column J width shrinks now that it no longer holds long mission text
# (closest value reachable through ColumnWidth's internal 1/6-character
# rounding to the recorded target width of 7.7265625)
$ws.Columns.Item(10).ColumnWidth = 6.8333333333333333

# --- Selection / scrolled view -----------------------------------------
$ws.Range("J3").Select()
